$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Reorder the goal rows (B3:F9) so that "Muss" (must-have) priority rows come
# first, followed by the "Soll" (should-have) rows - "Ziele der Prioritaet
# nach geordnet".
#
# Mapping of CURRENT row -> NEW row (content moves from key to value):
#   3 -> 5   4 -> 6   5 -> 7   6 -> 8   7 -> 9   8 -> 3   9 -> 4
#
# Strategy: stage the current rows 3-9 (values + formats) in a scratch area
# far below the used range, then paste them back into the sheet in the new
# order. This avoids clobbering source data while the destination ranges
# overlap the source ranges.
# ---------------------------------------------------------------------------

$srcRows = 3..9
$scratchBase = 300

foreach ($r in $srcRows) {
    $scratchRow = $scratchBase + $r
    $ws.Range("B$r`:F$r").Copy()
    $ws.Range("B$scratchRow`:F$scratchRow").PasteSpecial(-4122)
    $ws.Range("B$r`:F$r").Copy()
    $ws.Range("B$scratchRow`:F$scratchRow").PasteSpecial(-4163)
}
$excel.CutCopyMode = 0

$rowMap = @{3 = 5; 4 = 6; 5 = 7; 6 = 8; 7 = 9; 8 = 3; 9 = 4}

foreach ($r in $srcRows) {
    $scratchRow = $scratchBase + $r
    $destRow = $rowMap[$r]
    $ws.Range("B$scratchRow`:F$scratchRow").Copy()
    $ws.Range("B$destRow`:F$destRow").PasteSpecial(-4122)
    $ws.Range("B$scratchRow`:F$scratchRow").Copy()
    $ws.Range("B$destRow`:F$destRow").PasteSpecial(-4163)
}
$excel.CutCopyMode = 0

# clear the scratch area again
foreach ($r in $srcRows) {
    $scratchRow = $scratchBase + $r
    $ws.Range("B$scratchRow`:F$scratchRow").Clear()
}

# Column B ("Nr.") is a plain sequential row counter (1..7) independent of
# which content row it labels - restore it after the content shuffle above
# (which moved the old B values along with the rest of the row).
$n = 1
foreach ($r in $srcRows) {
    $ws.Range("B$r").Value = $n
    $n = $n + 1
}

# ---------------------------------------------------------------------------
# Re-apply left-alignment (in addition to the existing top alignment) across
# every data row - matches the style cleanup made alongside the reorder.
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("B3:F9")
$dataRange.VerticalAlignment = -4160     # xlTop
$dataRange.HorizontalAlignment = -4131   # xlLeft

# ---------------------------------------------------------------------------
# Row heights: the moved "Muss" rows (now rows 3 & 4) grow to fit their text
# at the top of the table; the rest keep their existing (travelling) height.
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 71.25
$ws.Rows.Item(4).RowHeight = 86.25
$ws.Rows.Item(5).RowHeight = 76.5
$ws.Rows.Item(6).RowHeight = 108.75
$ws.Rows.Item(7).RowHeight = 76.5
$ws.Rows.Item(8).RowHeight = 88
$ws.Rows.Item(9).RowHeight = 61.5
$ws.Rows.Item(10).RowHeight = 50.15
$ws.Rows.Item(11).RowHeight = 50.15

# ---------------------------------------------------------------------------
# Column widths: column D grows a bit wider than C/E/F.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 19.83
$ws.Columns.Item(4).ColumnWidth = 22.17
$ws.Columns.Item(5).ColumnWidth = 19.83
$ws.Columns.Item(6).ColumnWidth = 19.83

# ---------------------------------------------------------------------------
# Selection moves to K4 as part of the resave.
# ---------------------------------------------------------------------------
$ws.Range("K4").Select()
